$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: after "...is also pleasant to the eyes." add a new space-run,
# then a "_GoBack" bookmark wrapping a new sentence run.
# ---------------------------------------------------------------------------

# Merge the new sentence onto the existing run so the new text inherits the
# exact same run formatting (rStyle "eop", Arial, sz 28 / szCs 28, en-GB).
$rngEyes = $d.Content
$rngEyes.Find.Execute("is also pleasant to the eyes.", $true, $false, $false, $false, $false, $true, 1, $false, `
    "is also pleasant to the eyes. When they click the order button it will take them to a page where it clearly prints out success, and there will be a big cyan button they can click to get back. ", 2)

# Force a run split right after "eyes." (before the newly-added leading
# space) using a throwaway bookmark - inserting a bookmark always splits the
# enclosing run while both halves keep the original run properties.
$rngBoundary = $d.Content
$rngBoundary.Find.Execute("is also pleasant to the eyes.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boundaryPoint = $d.Range($rngBoundary.End, $rngBoundary.End)
$d.Bookmarks.Add("TempSplit1", $boundaryPoint)

# Wrap the new sentence in the real "_GoBack" bookmark (this also splits the
# run right before "When" so the leading-space run stays separate).
$rngSentence = $d.Content
$rngSentence.Find.Execute( `
    "When they click the order button it will take them to a page where it clearly prints out success, and there will be a big cyan button they can click to get back. ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("_GoBack", $rngSentence)

# Drop the throwaway bookmark now that it has done its job.
$d.Bookmarks("TempSplit1").Delete()

# ---------------------------------------------------------------------------
# Change 2: drop the stale <w:lastRenderedPageBreak/> before "Task 12:
# Refine the plan" (it is transient rendering metadata Word regenerates).
# A no-op text "replace" on that run is enough to make the host drop it.
# ---------------------------------------------------------------------------
$rngTask12 = $d.Content
$rngTask12.Find.Execute("Task 12: Refine the plan", $true, $false, $false, $false, $false, $true, 1, $false, "Task 12: Refine the plan", 2)

# ---------------------------------------------------------------------------
# Change 3: merge " I am happy the way it has turned out," and " and am
# ready for version 4.0" into a single run. The old "_GoBack" bookmark that
# used to sit between them was already relocated to its new spot above (a
# document can only have one bookmark per name, so Bookmarks.Add in Change 1
# already pulled the "_GoBack" name off of this old location).
# ---------------------------------------------------------------------------

# Keep " did the steps very easily." from merging into the run(s) after it
# by planting another throwaway bookmark right at that boundary.
$rngEasily = $d.Content
$rngEasily.Find.Execute(" did the steps very easily.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boundaryPoint2 = $d.Range($rngEasily.End, $rngEasily.End)
$d.Bookmarks.Add("TempSplit2", $boundaryPoint2)

# Belt-and-braces: if for any reason a stray "_GoBack" bookmark still sits
# at this old location (after the range we just protected), drop it so it
# doesn't block the run merge below.
if ($d.Bookmarks.Exists("_GoBack") -and ($d.Bookmarks("_GoBack").Start -gt $rngEasily.End)) {
    $d.Bookmarks("_GoBack").Delete()
}

# Touch the text so the two identically-formatted (plain) runs coalesce.
$rngHappy = $d.Content
$rngHappy.Find.Execute("turned out,", $true, $false, $false, $false, $false, $true, 1, $false, "turned out,", 2)

# Drop the throwaway bookmark.
$d.Bookmarks("TempSplit2").Delete()
